$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header row (e.g. H1: bold, bordered, centered).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New data columns I (I0) and J (IF) for rows 2-29
$data = @{
  2  = @(6, 9)
  3  = @(5, 8)
  4  = @(5, 8)
  5  = @(2, 7)
  6  = @(2, 6)
  7  = @(1, 6)
  8  = @(2, 7)
  9  = @(1, 7)
  10 = @(1, 6)
  11 = @(1, 5)
  12 = @(1, 4)
  13 = @(1, 2)
  14 = @(6, 8)
  15 = @(1, 4)
  16 = @(2, 5)
  17 = @(1, 5)
  18 = @(1, 5)
  19 = @(1, 4)
  20 = @(6, 8)
  21 = @(6, 9)
  22 = @(7, 8)
  23 = @(1, 5)
  24 = @(1, 6)
  25 = @(1, 5)
  26 = @(1, 5)
  27 = @(1, 4)
  28 = @(4, 6)
  29 = @(3, 4)
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  $ws.Cells.Item($row, 9).Value = $vals[0]
  $ws.Cells.Item($row, 10).Value = $vals[1]
}
